$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update win/loss probability matrix cells with refreshed simulation values
# (added more games, sped up simulate game logic, drafted optimization logic)

# Row 2
$ws.Range("B2").Value = 0.204
$ws.Range("C2").Value = 0.524
$ws.Range("J2").Value = 0.028
$ws.Range("P2").Value = 0.164
$ws.Range("S2").Value = 0.08

# Row 3
$ws.Range("C3").Value = 0.02898550724637681
$ws.Range("J3").Value = 0.01449275362318841
$ws.Range("P3").Value = 0.7536231884057971
$ws.Range("S3").Value = 0.2028985507246377

# Row 4
$ws.Range("J4").Value = 0.025
$ws.Range("O4").Value = 0.025
$ws.Range("P4").Value = 0.65
$ws.Range("S4").Value = 0.3

# Row 6
$ws.Range("B6").Value = 0.03896103896103896
$ws.Range("D6").Value = 0.01298701298701299
$ws.Range("F6").Value = 0.06926406926406926
$ws.Range("J6").Value = 0.2640692640692641
$ws.Range("O6").Value = 0.03463203463203463
$ws.Range("Q6").Value = 0.1515151515151515
$ws.Range("R6").Value = 0.06493506493506493
$ws.Range("S6").Value = 0.3636363636363636

# Row 7
$ws.Range("B7").Value = 0.1012658227848101
$ws.Range("D7").Value = 0.02531645569620253
$ws.Range("F7").Value = 0.0949367088607595
$ws.Range("J7").Value = 0.1392405063291139
$ws.Range("O7").Value = 0.01265822784810127
$ws.Range("Q7").Value = 0.120253164556962
$ws.Range("R7").Value = 0.0949367088607595
$ws.Range("S7").Value = 0.4113924050632912

# Row 8
$ws.Range("B8").Value = 0.07142857142857142
$ws.Range("D8").Value = 0.01612903225806452
$ws.Range("F8").Value = 0.06912442396313365
$ws.Range("J8").Value = 0.1589861751152074
$ws.Range("O8").Value = 0.02995391705069124
$ws.Range("Q8").Value = 0.1774193548387097
$ws.Range("R8").Value = 0.09216589861751152
$ws.Range("S8").Value = 0.3847926267281106

# Row 9
$ws.Range("B9").Value = 0.08717948717948718
$ws.Range("D9").Value = 0.01025641025641026
$ws.Range("F9").Value = 0.07692307692307693
$ws.Range("J9").Value = 0.1282051282051282
$ws.Range("O9").Value = 0.03589743589743589
$ws.Range("Q9").Value = 0.1333333333333333
$ws.Range("R9").Value = 0.09230769230769231
$ws.Range("S9").Value = 0.4358974358974359

# Row 10
$ws.Range("B10").Value = 0.1019900497512438
$ws.Range("D10").Value = 0.02155887230514096
$ws.Range("E10").Value = 0.001658374792703151
$ws.Range("F10").Value = 0.07545605306799337
$ws.Range("J10").Value = 0.1434494195688226
$ws.Range("O10").Value = 0.02321724709784411
$ws.Range("Q10").Value = 0.1981757877280265
$ws.Range("R10").Value = 0.0812603648424544
$ws.Range("S10").Value = 0.3532338308457711

# Row 11
$ws.Range("G11").Value = 0.1062992125984252
$ws.Range("J11").Value = 0.1102362204724409
$ws.Range("K11").Value = 0.1811023622047244
$ws.Range("L11").Value = 0.5905511811023622
$ws.Range("S11").Value = 0.01181102362204724

# Row 12
$ws.Range("G12").Value = 0.7483443708609272
$ws.Range("J12").Value = 0.2052980132450331
$ws.Range("K12").Value = 0.006622516556291391
$ws.Range("L12").Value = 0.02649006622516556
$ws.Range("S12").Value = 0.01324503311258278

# Row 13
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.275
$ws.Range("S13").Value = 0.1

# Row 15
$ws.Range("F15").Value = 0.01304347826086956
$ws.Range("H15").Value = 0.1652173913043478
$ws.Range("I15").Value = 0.08260869565217391
$ws.Range("J15").Value = 0.308695652173913
$ws.Range("K15").Value = 0.06086956521739131
$ws.Range("M15").Value = 0.01304347826086956
$ws.Range("O15").Value = 0.04347826086956522
$ws.Range("S15").Value = 0.3130434782608696

# Row 16
$ws.Range("F16").Value = 0.01204819277108434
$ws.Range("H16").Value = 0.2349397590361446
$ws.Range("I16").Value = 0.09036144578313253
$ws.Range("J16").Value = 0.4096385542168675
$ws.Range("K16").Value = 0.06024096385542169
$ws.Range("M16").Value = 0.03012048192771084
$ws.Range("O16").Value = 0.04819277108433735
$ws.Range("S16").Value = 0.1144578313253012

# Row 17
$ws.Range("F17").Value = 0.01526717557251908
$ws.Range("H17").Value = 0.1908396946564886
$ws.Range("I17").Value = 0.1119592875318066
$ws.Range("J17").Value = 0.3994910941475827
$ws.Range("K17").Value = 0.09923664122137404
$ws.Range("M17").Value = 0.02290076335877863
$ws.Range("O17").Value = 0.07888040712468193
$ws.Range("S17").Value = 0.08142493638676845

# Row 18
$ws.Range("F18").Value = 0.03243243243243243
$ws.Range("H18").Value = 0.2108108108108108
$ws.Range("I18").Value = 0.08108108108108109
$ws.Range("J18").Value = 0.4594594594594595
$ws.Range("K18").Value = 0.07027027027027027
$ws.Range("M18").Value = 0.01081081081081081
$ws.Range("O18").Value = 0.05405405405405406
$ws.Range("S18").Value = 0.08108108108108109

# Row 19
$ws.Range("F19").Value = 0.01858407079646018
$ws.Range("H19").Value = 0.2194690265486726
$ws.Range("I19").Value = 0.09203539823008849
$ws.Range("J19").Value = 0.3610619469026549
$ws.Range("K19").Value = 0.1176991150442478
$ws.Range("M19").Value = 0.02035398230088496
$ws.Range("S19").Value = 0.09823008849557523
